$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update benchmark values (formulas in column D recalc automatically)
$ws.Range("C27").Value = 33723.286738000003
$ws.Range("C28").Value = 17065.508472000001

# Update the active cell selection to match the saved state
$ws.Range("C27").Select()
